# Add emulated x64 benchmark results for Qualcomm Snapdragon 835.
# Inserts a new data row at row 16 (existing rows 16-23 shift to 17-24)
# and fills it in with the new benchmark numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 16, pushing everything else down.
$ws.Rows.Item(16).Insert() | Out-Null
$ws.Rows.Item(16).Select() | Out-Null

$ws.Cells.Item(16, 1).Value = "Qualcomm"
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(16, 3).Value = "Snapdragon 835"
$ws.Cells.Item(16, 4).Value = 3.5
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 2.2
$ws.Cells.Item(16, 8).Value = 2.4
$ws.Cells.Item(16, 9).Value = "x86-64 (emu)"
$ws.Cells.Item(16, 10).Value = "Balanced"
$ws.Cells.Item(16, 11).Value = 6
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = "LPDDR4x"
$ws.Cells.Item(16, 14).Value = 1866
$ws.Cells.Item(16, 15).Value = 2.41
$ws.Cells.Item(16, 16).Value = 4.97
$ws.Cells.Item(16, 17).Value = 10.05
$ws.Cells.Item(16, 18).Value = 21.07
